$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the sanity-system tasks ("Create sprites for columns", "Compose
# attack SFX", "Compose player SFX", "Create sprite for beloved") as
# Finished ("Yes") in the "Finished?" column (E). E25 had no value yet;
# E28/E29/E32 were previously "WIP".
$ws.Range("E25").Value = "Yes"
$ws.Range("E28").Value = "Yes"
$ws.Range("E29").Value = "Yes"
$ws.Range("E32").Value = "Yes"

# Leave the view scrolled/selected where the edits were made.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E26").Select()
